$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new comment text for PRICE in column E, row 3 (matches the existing
# E2 comment pattern for TMTID/EFFECTIVE_DATE)
$ws.Range("E3").Value = "PRICE ไม่ต้องใส่หน่วย"

# Update the selected cell / sqref in the sheet view from G7 to G11
$ws.Range("G11").Select()
